$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Row 3: "sendkeys / upload image" row becomes the "click" row (only A & B survive,
# C keeps the hyperlink style but no value, D is cleared)
$ws.Cells.Item(3, 1).Value = "click"
$ws.Cells.Item(3, 3).Value = ""
$ws.Cells.Item(3, 3).Style = "Hyperlink"
$ws.Cells.Item(3, 4).Value = ""

# Row 4: new "exe / Autoit.exe / upload image" row
$ws.Cells.Item(4, 1).Value = "exe"
$ws.Cells.Item(4, 2).Value = ""
$ws.Cells.Item(4, 3).Value = "C:\Users\jack.zhong\Desktop\Autoit.exe"
$ws.Cells.Item(4, 4).Value = "上传图片"

# Row 5: "text" row now points at label[1] and uses the new A / 萨达 values
$ws.Cells.Item(5, 2).Value = ".//*[@id='radio']/label[1]"
$ws.Cells.Item(5, 3).Value = "A"
$ws.Cells.Item(5, 4).Value = "萨达"

# Row 6: becomes the "contain" row
$ws.Cells.Item(6, 1).Value = "contain"
$ws.Cells.Item(6, 2).Value = ".//*[@id='radio']/label[1]"
$ws.Cells.Item(6, 3).Value = "A>百"
$ws.Cells.Item(6, 4).Value = "duany"

# Row 7: brand new "wait / [1]" row
$ws.Cells.Item(7, 1).Value = "wait"
$ws.Cells.Item(7, 3).Value = "[1]"

# Selection moves from A6 to C6
$ws.Range("C6").Select()

# Three blank rows were inserted above the old row 11 (its 15.75pt custom height
# slides down to row 14, and the sheet grows from 15 to 18 rows)
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(11).Insert()
